$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 549
$ws.Range("F4").Value = 5810
$ws.Range("F9").Value = 1558
$ws.Range("F11").Value = 26
$ws.Range("F12").Value = 657
$ws.Range("F13").Value = 1571
$ws.Range("F14").Value = 1500
$ws.Range("F16").Value = 110
$ws.Range("F17").Value = 595
$ws.Range("F18").Value = 4275
$ws.Range("F19").Value = 10
$ws.Range("F22").Value = 798
$ws.Range("F23").Value = 1
$ws.Range("F25").Value = 15
$ws.Range("F26").Value = 2266
$ws.Range("F29").Value = 11
$ws.Range("F30").Value = 446
$ws.Range("F31").Value = 1211
$ws.Range("F32").Value = 778
$ws.Range("F34").Value = 1155
$ws.Range("F35").Value = 1164

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 104
$ws.Range("F15").Value = 17
$ws.Range("F18").Value = 119
$ws.Range("F19").Value = 286
$ws.Range("F20").Value = 223
$ws.Range("F21").Value = 487
$ws.Range("F23").Value = 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 245
$ws.Range("F3").Value = 598
$ws.Range("F4").Value = 142
$ws.Range("F5").Value = 239

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 245
$ws.Range("F3").Value = 549
$ws.Range("F6").Value = 598
$ws.Range("F7").Value = 142
$ws.Range("F8").Value = 5810
$ws.Range("F14").Value = 104
$ws.Range("F20").Value = 1558
$ws.Range("F23").Value = 1571
$ws.Range("F25").Value = 1500
$ws.Range("F27").Value = 110
$ws.Range("F28").Value = 595
$ws.Range("F30").Value = 4275
$ws.Range("F33").Value = 798
$ws.Range("F34").Value = 2266
$ws.Range("F37").Value = 446
$ws.Range("F38").Value = 1211
$ws.Range("F40").Value = 119
$ws.Range("F41").Value = 286
$ws.Range("F42").Value = 223
$ws.Range("F43").Value = 487
$ws.Range("F44").Value = 778
$ws.Range("F46").Value = 1155
$ws.Range("F47").Value = 2
$ws.Range("F48").Value = 1164
